$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Raise the price of 1 pizza from 4 to 40 (Enter price of 1 pizza)
$ws.Range("F6").Value = 40

# Force a full recalculation so dependent formulas, cached chart values,
# and the "optimal pizza/book" lookups refresh for the new price.
$excel.CalculateFullRebuild()

# Match number formatting of N13 (Opt. book) to N12 (Opt. pizza) - both
# should show as 0.00 now that N13's offset result changed.
$ws.Range("N13").NumberFormat = $ws.Range("N12").NumberFormat
